# Apply "put some check marks, which functions is needed to work on" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the TODO/status notes in column A (rows 22, 24, 25, 27, 28, 29)
$ws.Range("A22").Value = "1x 116 check?"
$ws.Range("A23").Value = "2x 120 check ? "
$ws.Range("A24").Value = "1x 122 check?"
$ws.Range("A25").Value = "2x 123 check? Basic (working on it"
$ws.Range("A26").Value = "2x 127 check? "
$ws.Range("A27").Value = "1x 128 check? Working on damage"
$ws.Range("A28").Value = "2x 134 check? Function needed for flipCoin"
$ws.Range("A29").Value = "1x 135 check?"

# Widen column A to fit the new, longer text (~56.83 "characters", as the
# workbook XML records after Excel's pixel-rounding of the column width)
$ws.Columns.Item(1).ColumnWidth = 56

# Move selection to A29, matching the author's last-saved cursor position
$ws.Range("A29").Select()
